$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4274099.5
$ws.Range("I19").Value = 6945163
$ws.Range("K19").Value = 6945163
$ws.Range("M19").Value = -6944988
$ws.Range("H33").Value = 1840.1177
$ws.Range("I33").Value = 436.5
$ws.Range("J33").Value = 2605.7273
$ws.Range("K33").Value = 436.5
$ws.Range("L33").Value = 2605.7273
$ws.Range("M33").Value = -207.5
$ws.Range("N33").Value = -3063.7273
$ws.Range("H98").Value = 9763.5
$ws.Range("I98").Value = 12235.777
$ws.Range("K98").Value = 12235.777
$ws.Range("M98").Value = -10737.777
$ws.Range("H122").Value = 9763.5
$ws.Range("I122").Value = 12235.777
$ws.Range("K122").Value = 36707.331
$ws.Range("M122").Value = -34257.331
$ws.Range("H135").Value = 3171.4443
$ws.Range("I135").Value = 3171.4443
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 28542.9987
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -26007.9987
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 7417988.5
$ws.Range("I138").Value = 13079.889
$ws.Range("J138").Value = 18525352
$ws.Range("K138").Value = 39239.667
$ws.Range("L138").Value = 55576056
$ws.Range("M138").Value = -34099.667
$ws.Range("N138").Value = -55586336
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9278.733
$ws.Range("I45").Value = 11089.454
$ws.Range("K45").Value = 11089.454
$ws.Range("M45").Value = -10712.454
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 26324590
$ws.Range("I86").Value = 8574.5
$ws.Range("J86").Value = 166676670
$ws.Range("K86").Value = 8574.5
$ws.Range("L86").Value = 166676670
$ws.Range("M86").Value = -7451.5
$ws.Range("N86").Value = -166678916
$ws.Range("H89").Value = 26324590
$ws.Range("I89").Value = 8574.5
$ws.Range("J89").Value = 166676670
$ws.Range("K89").Value = 42872.5
$ws.Range("L89").Value = 833383350
$ws.Range("M89").Value = -37256.5
$ws.Range("N89").Value = -833394582
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3585.0833
$ws.Range("I31").Value = 2036.2
$ws.Range("K31").Value = 2036.2
$ws.Range("M31").Value = -1741.2
$ws.Range("H34").Value = 3585.0833
$ws.Range("I34").Value = 2036.2
$ws.Range("K34").Value = 2036.2
$ws.Range("M34").Value = -1834.2
$ws.Range("H51").Value = 40000
$ws.Range("J51").Value = 40000
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41472
$ws.Range("H61").Value = 40000
$ws.Range("J61").Value = 40000
$ws.Range("L61").Value = 40000
$ws.Range("N61").Value = -40696
$ws.Range("H99").Value = 2489
$ws.Range("I99").Value = 2489
$ws.Range("K99").Value = 2489
$ws.Range("M99").Value = -991
$ws.Range("H107").Value = 675.4706
$ws.Range("I107").Value = 675.4706
$ws.Range("K107").Value = 675.4706
$ws.Range("M107").Value = 1244.5294
$ws.Range("H126").Value = 2489
$ws.Range("I126").Value = 2489
$ws.Range("K126").Value = 7467
$ws.Range("M126").Value = -4997
$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 60000
$ws.Range("L127").Value = 60000
$ws.Range("N127").Value = -69920
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 165.63637
$ws.Range("I12").Value = 133.75
$ws.Range("J12").Value = 183.85715
$ws.Range("K12").Value = 401.25
$ws.Range("L12").Value = 551.5714499999999
$ws.Range("M12").Value = -228.25
$ws.Range("N12").Value = -897.5714499999999
$ws.Range("H131").Value = 3337638
$ws.Range("I131").Value = 2000
$ws.Range("J131").Value = 3513197.8
$ws.Range("K131").Value = 6000
$ws.Range("L131").Value = 10539593.4
$ws.Range("M131").Value = -960
$ws.Range("N131").Value = -10549673.4
$ws.Range("H132").Value = 2691.6155
$ws.Range("J132").Value = 2817.818
$ws.Range("L132").Value = 25360.362
$ws.Range("N132").Value = -30420.362
$ws.Range("H140").Value = 3729.9092
$ws.Range("I140").Value = 2878.625
$ws.Range("K140").Value = 8635.875
$ws.Range("M140").Value = -3455.875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8784.857
$ws.Range("I43").Value = 3000
$ws.Range("K43").Value = 3000
$ws.Range("M43").Value = -2849
$ws.Range("H126").Value = 2471.25
$ws.Range("I126").Value = 1774
$ws.Range("K126").Value = 5322
$ws.Range("M126").Value = -2852
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1117.6111
$ws.Range("I46").Value = 734.2
$ws.Range("J46").Value = 1265.0769
$ws.Range("K46").Value = 734.2
$ws.Range("L46").Value = 1265.0769
$ws.Range("M46").Value = -546.2
$ws.Range("N46").Value = -1641.0769
$ws.Range("H68").Value = 6972.5
$ws.Range("I68").Value = 4999.6
$ws.Range("K68").Value = 4999.6
$ws.Range("M68").Value = -4250.6
$ws.Range("H71").Value = 6972.5
$ws.Range("I71").Value = 4999.6
$ws.Range("K71").Value = 24998
$ws.Range("M71").Value = -21254
$ws.Range("H82").Value = 2871.4285
$ws.Range("I82").Value = 500
$ws.Range("J82").Value = 3266.6667
$ws.Range("K82").Value = 500
$ws.Range("L82").Value = 3266.6667
$ws.Range("M82").Value = -139
$ws.Range("N82").Value = -3988.6667
$ws.Range("H85").Value = 2871.4285
$ws.Range("I85").Value = 500
$ws.Range("J85").Value = 3266.6667
$ws.Range("K85").Value = 500
$ws.Range("L85").Value = 3266.6667
$ws.Range("M85").Value = 748
$ws.Range("N85").Value = -5762.6667
$ws.Range("H101").Value = 10517.625
$ws.Range("J101").Value = 10517.625
$ws.Range("L101").Value = 10517.625
$ws.Range("N101").Value = -17007.625
$ws.Range("H132").Value = 3923.4285
$ws.Range("I132").Value = 3558.7
$ws.Range("J132").Value = 4255
$ws.Range("K132").Value = 10676.1
$ws.Range("L132").Value = 12765
$ws.Range("M132").Value = -8146.099999999999
$ws.Range("N132").Value = -17825
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1082.6842
$ws.Range("I113").Value = 706.0909
$ws.Range("J113").Value = 1600.5
$ws.Range("K113").Value = 2118.2727
$ws.Range("L113").Value = 4801.5
$ws.Range("M113").Value = 51.72730000000001
$ws.Range("N113").Value = -9141.5
$ws.Range("H122").Value = 3186.842
$ws.Range("I122").Value = 1955.2916
$ws.Range("K122").Value = 5865.8748
$ws.Range("M122").Value = -3415.8748

Write-Host "Applied all market-price updates."